$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp text in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 22:22"

# --- Update province data rows (Ciudad, Casos totales, Casos activos, Recuperados, Muertes) ---
# Asturias overtook Sevilla and Gipuzkoa/Guipuzcoa in total cases, so the three
# rows re-sort (labels move, figures move with them); Segovia stays put.
$ws.Range("A21").Value = "Asturias"
$ws.Range("B21").Value = 2348
$ws.Range("C21").Value = 599
$ws.Range("D21").Value = 1549
$ws.Range("E21").Value = 200

$ws.Range("A22").Value = "Sevilla"
$ws.Range("B22").Value = 2329
$ws.Range("C22").Value = 459
$ws.Range("D22").Value = 1658
$ws.Range("E22").Value = 212

$ws.Range("A23").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B23").Value = 2328
$ws.Range("C23").Value = 6144
$ws.Range("D23").Value = 4953
$ws.Range("E23").Value = 209

# Murcia (row 30)
$ws.Range("B30").Value = 1659
$ws.Range("C30").Value = 662
$ws.Range("D30").Value = 880
$ws.Range("E30").Value = 117

# Tenerife (row 34)
$ws.Range("C34").Value = 813
$ws.Range("E34").Value = 119

# Gran Canaria (row 50)
$ws.Range("B50").Value = 502
$ws.Range("D50").Value = 234

# Ceuta (row 54)
$ws.Range("B54").Value = 115
$ws.Range("C54").Value = 73
$ws.Range("D54").Value = 38

# Melilla (row 55)
$ws.Range("B55").Value = 109
$ws.Range("C55").Value = 46
$ws.Range("D55").Value = 61
